# The document has three logo images living in the headers/footers:
#   - Header 2  : "BTec_Logo-Orange"  (image1.jpg -> image2.jpg)
#   - Footer 1  : "PearsonLogo.png"   (image2.png -> image1.png)
#   - Footer 2  : "PearsonLogo.png"   (image2.png -> image1.png)
# Only the image's friendly "name" metadata changes (descr/ids/embeds
# are untouched) - this mirrors Word silently renumbering the inline
# pictures when the document was re-saved.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Header 2: BTec logo -> image2.jpg -------------------------------
$hdr2 = $sec.Headers.Item(2)
if ($hdr2.Exists -and $hdr2.Range.InlineShapes.Count -gt 0) {
    $btecShape = $hdr2.Range.InlineShapes.Item(1)
    $btecShape.Name = "image2.jpg"
}

# --- Footer 1: Pearson logo -> image1.png -----------------------------
$ftr1 = $sec.Footers.Item(1)
if ($ftr1.Exists -and $ftr1.Range.InlineShapes.Count -gt 0) {
    $pearsonShape1 = $ftr1.Range.InlineShapes.Item(1)
    $pearsonShape1.Name = "image1.png"
}

# --- Footer 2: Pearson logo -> image1.png -----------------------------
$ftr2 = $sec.Footers.Item(2)
if ($ftr2.Exists -and $ftr2.Range.InlineShapes.Count -gt 0) {
    $pearsonShape2 = $ftr2.Range.InlineShapes.Item(1)
    $pearsonShape2.Name = "image1.png"
}
